$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 750.4737
$ws.Range("I33").Value = 686.1177
$ws.Range("J33").Value = 1297.5
$ws.Range("K33").Value = 686.1177
$ws.Range("L33").Value = 1297.5
$ws.Range("M33").Value = -457.1177
$ws.Range("N33").Value = -1755.5
$ws.Range("H94").Value = 4690.143
$ws.Range("I94").Value = 2138.5
$ws.Range("K94").Value = 2138.5
$ws.Range("M94").Value = -1687.5
$ws.Range("H129").Value = 2443.3684
$ws.Range("I129").Value = 2332.6155
$ws.Range("K129").Value = 6997.8465
$ws.Range("M129").Value = -1997.8465
$ws.Range("H137").Value = 2498
$ws.Range("I137").Value = 1561.2632
$ws.Range("K137").Value = 4683.7896
$ws.Range("M137").Value = -2133.7896
$ws.Range("H140").Value = 78431
$ws.Range("J140").Value = 78431
$ws.Range("L140").Value = 78431
$ws.Range("N140").Value = -88791

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4340.9165
$ws.Range("I45").Value = 4008.2727
$ws.Range("K45").Value = 4008.2727
$ws.Range("M45").Value = -3631.2727
$ws.Range("H61").Value = 5710.6924
$ws.Range("I61").Value = 3260.5386
$ws.Range("K61").Value = 3260.5386
$ws.Range("M61").Value = -3048.5386
$ws.Range("H74").Value = 2362767.8
$ws.Range("I74").Value = 2979771.8
$ws.Range("J74").Value = 6934.273
$ws.Range("K74").Value = 2979771.8
$ws.Range("L74").Value = 6934.273
$ws.Range("M74").Value = -2978897.8
$ws.Range("N74").Value = -8682.273000000001
$ws.Range("H77").Value = 2362767.8
$ws.Range("I77").Value = 2979771.8
$ws.Range("J77").Value = 6934.273
$ws.Range("K77").Value = 14898859
$ws.Range("L77").Value = 34671.365
$ws.Range("M77").Value = -14894491
$ws.Range("N77").Value = -43407.365
$ws.Range("H122").Value = 1625.9487
$ws.Range("I122").Value = 1663.1945
$ws.Range("K122").Value = 4989.583500000001
$ws.Range("M122").Value = -2539.583500000001
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840
$ws.Range("H132").Value = 479198.6
$ws.Range("I132").Value = 584911
$ws.Range("J132").Value = 9365.777
$ws.Range("K132").Value = 1754733
$ws.Range("L132").Value = 28097.331
$ws.Range("M132").Value = -1752203
$ws.Range("N132").Value = -33157.331
$ws.Range("H136").Value = 5710.6924
$ws.Range("I136").Value = 3260.5386
$ws.Range("K136").Value = 9781.6158
$ws.Range("M136").Value = -7231.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 701.2778
$ws.Range("I94").Value = 576.9524
$ws.Range("J94").Value = 875.3333
$ws.Range("K94").Value = 576.9524
$ws.Range("L94").Value = 875.3333
$ws.Range("M94").Value = -125.9524
$ws.Range("N94").Value = -1777.3333
$ws.Range("H105").Value = 3193.1
$ws.Range("I105").Value = 3118.9412
$ws.Range("K105").Value = 3118.9412
$ws.Range("M105").Value = -1371.9412
$ws.Range("H134").Value = 378858.78
$ws.Range("I134").Value = 437438.03
$ws.Range("J134").Value = 7857
$ws.Range("K134").Value = 1312314.09
$ws.Range("L134").Value = 23571
$ws.Range("M134").Value = -1309779.09
$ws.Range("N134").Value = -28641

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10581.434
$ws.Range("I31").Value = 2250
$ws.Range("K31").Value = 2250
$ws.Range("M31").Value = -1955
$ws.Range("H34").Value = 10581.434
$ws.Range("I34").Value = 2250
$ws.Range("K34").Value = 2250
$ws.Range("M34").Value = -2048
$ws.Range("H122").Value = 1809.8572
$ws.Range("I122").Value = 1822
$ws.Range("J122").Value = 1779.5
$ws.Range("K122").Value = 5466
$ws.Range("L122").Value = 5338.5
$ws.Range("M122").Value = -3016
$ws.Range("N122").Value = -10238.5
$ws.Range("H132").Value = 6421837
$ws.Range("I132").Value = 12051.484
$ws.Range("J132").Value = 41675656
$ws.Range("K132").Value = 36154.452
$ws.Range("L132").Value = 125026968
$ws.Range("M132").Value = -33624.452
$ws.Range("N132").Value = -125032028
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060
$ws.Range("H139").Value = 83333
$ws.Range("J139").Value = 83333
$ws.Range("L139").Value = 83333
$ws.Range("N139").Value = -93613

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 838.2727
$ws.Range("J107").Value = 945.4286
$ws.Range("L107").Value = 2836.2858
$ws.Range("N107").Value = -6676.2858
$ws.Range("H120").Value = 21343.334
$ws.Range("I120").Value = 14498.5
$ws.Range("J120").Value = 35033
$ws.Range("K120").Value = 43495.5
$ws.Range("L120").Value = 105099
$ws.Range("M120").Value = -38657.5
$ws.Range("N120").Value = -114775

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2600.9697
$ws.Range("I132").Value = 1855.8276
$ws.Range("J132").Value = 8003.25
$ws.Range("K132").Value = 5567.4828
$ws.Range("L132").Value = 24009.75
$ws.Range("M132").Value = -3037.4828
$ws.Range("N132").Value = -29069.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 66552.42999999999
$ws.Range("I74").Value = 62499.5
$ws.Range("J74").Value = 68173.60000000001
$ws.Range("K74").Value = 62499.5
$ws.Range("L74").Value = 68173.60000000001
$ws.Range("M74").Value = -61501.5
$ws.Range("N74").Value = -70169.60000000001
$ws.Range("H77").Value = 66552.42999999999
$ws.Range("I77").Value = 62499.5
$ws.Range("J77").Value = 68173.60000000001
$ws.Range("K77").Value = 187498.5
$ws.Range("L77").Value = 204520.8
$ws.Range("M77").Value = -182506.5
$ws.Range("N77").Value = -214504.8
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H136").Value = 8033.0884
$ws.Range("I136").Value = 7649.2256
$ws.Range("K136").Value = 22947.6768
$ws.Range("M136").Value = -20397.6768

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2398.9285
$ws.Range("I122").Value = 2179.6155
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 6538.8465
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -4088.8465
$ws.Range("N122").Value = -20650
$ws.Range("H136").Value = 12698969
$ws.Range("I136").Value = 15237579
$ws.Range("K136").Value = 45712737
$ws.Range("M136").Value = -45710187
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 65999.60000000001
$ws.Range("J141").Value = 74999.5
$ws.Range("L141").Value = 74999.5
$ws.Range("N141").Value = -85359.5
